$d = $word.ActiveDocument

# --- Step 1: find the last occurrence of the target URL paragraph and
#     turn it into a hyperlink styled with the document's "Hyperlink"
#     character style (maps to styleId "ac" in this document). ---
$targetUrl = "https://www.data.go.kr/tcs/eds/selectCoreDataView.do"

$targetPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "$targetUrl*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "target paragraph not found"
}

$linkRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
$d.Hyperlinks.Add($linkRange, $targetUrl)

# Re-acquire the run range (exclude the paragraph mark) and normalize its
# character style to the document's Hyperlink style (styleId "ac").
$afterLinkPara = $targetPara
$styleRange = $d.Range($afterLinkPara.Range.Start, $afterLinkPara.Range.End - 1)
$styleRange.Style = "Hyperlink"

# --- Step 2: replace the final (empty, sz=16) paragraph with a run of
#     five fresh paragraphs: a blank one, then four plain-text Korean
#     labels, each using rFonts hint=eastAsia and no inherited pPr. ---
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$last = $d.Paragraphs.Item($d.Paragraphs.Count)

$xml = "<w:p $ns/>" +
       "<w:p $ns><w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>통계(Rank)</w:t></w:r></w:p>" +
       "<w:p $ns><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>내보내기</w:t></w:r></w:p>" +
       "<w:p $ns><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>시각화</w:t></w:r></w:p>" +
       "<w:p $ns><w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>조회</w:t></w:r></w:p>"

$last.Range.InsertXML($xml)

Write-Output "done"
